$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add new row 9 data, copying formatting from row 8 first (values set before formats
#    so numeric cells don't get coerced into text by a text-like number format).
$ws.Range("A9").Value = 44751
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = "PILLACA RIVERA"
$ws.Range("D9").Value = "CANDY VANESSA"
$ws.Range("E9").Value = 73908404
$ws.Range("F9").Value = "VACACIONES"
$ws.Range("G9").Value = "Pago incompleto"
$ws.Range("H9").Value = "Falto pagar 6 dias regulares"

$src = $ws.Range("A8:H8")
$dst = $ws.Range("A9:H9")
$src.Copy()
$dst.PasteSpecial(-4122)

# 2. Rebuild the AutoFilter over the new range, filtering column C (index 3) on the new
#    last name, and drop the old date filterColumn.
$ws.AutoFilterMode = $false
$crit = @("PILLACA RIVERA")
[void]$ws.Range("A2:H9").AutoFilter(3, $crit, 7)

# 3. Row 8 should now be hidden (it no longer matches the active filter).
$ws.Range("A8").EntireRow.Hidden = $true

# 4. Update the active selection shown when the sheet is opened.
[void]$ws.Range("H11").Select()

# 5. Keep the hidden _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Hoja1!_FilterDatabase") {
        $n.RefersTo = "=Hoja1!`$A`$2:`$H`$9"
    }
}
